# Apply the edit described by the diff:
#  - Two new weekly records are inserted right before the existing row 236
#    (pushing the old rows 236..324 down to 238..326), so the sheet's
#    dimension grows from A1:R324 to A1:R326.
#  - The two newly inserted rows (236 and 237) are populated with their
#    own data (new dates / volumes / prices), everything else on the
#    sheet stays exactly as it was (just shifted down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above what is currently row 236. This shifts the
# old row 236 down to row 238, old row 237 to row 239, and so on, all the
# way down to old row 324 becoming row 326 - matching the new dimension.
$ws.Rows.Item(236).Insert()
$ws.Rows.Item(236).Insert()

# Common (constant) columns shared by every data row in this sheet.
$mercadoId = 7
$mercado = "Terminal Hortofrutícola Agro Chillán"
$region = "Ñuble"
$codreg = 16
$categoriaId = 100112009
$categoria = "Acelga"
$variedad = "Sin especificar"
$origen = "Provincia de Diguillín"
$clasificacion = "Hortaliza"

# New row 236: Primera, $/atado 0,5 a 1 kilo
$ws.Cells.Item(236, 1).Value = $mercadoId
$ws.Cells.Item(236, 2).Value = $mercado
$ws.Cells.Item(236, 3).Value = $region
$ws.Cells.Item(236, 4).Value = 44900
$ws.Cells.Item(236, 5).Value = $codreg
$ws.Cells.Item(236, 6).Value = $categoriaId
$ws.Cells.Item(236, 7).Value = $categoria
$ws.Cells.Item(236, 8).Value = $variedad
$ws.Cells.Item(236, 9).Value = "Primera"
$ws.Cells.Item(236, 10).Value = 400
$ws.Cells.Item(236, 11).Value = 600
$ws.Cells.Item(236, 12).Value = 700
$ws.Cells.Item(236, 13).Value = 650
$ws.Cells.Item(236, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(236, 15).Value = $origen
$ws.Cells.Item(236, 16).Value = 650
$ws.Cells.Item(236, 17).Value = 1
$ws.Cells.Item(236, 18).Value = $clasificacion

# New row 237: Segunda, $/atado 0,5 a 1 kilo
$ws.Cells.Item(237, 1).Value = $mercadoId
$ws.Cells.Item(237, 2).Value = $mercado
$ws.Cells.Item(237, 3).Value = $region
$ws.Cells.Item(237, 4).Value = 44900
$ws.Cells.Item(237, 5).Value = $codreg
$ws.Cells.Item(237, 6).Value = $categoriaId
$ws.Cells.Item(237, 7).Value = $categoria
$ws.Cells.Item(237, 8).Value = $variedad
$ws.Cells.Item(237, 9).Value = "Segunda"
$ws.Cells.Item(237, 10).Value = 300
$ws.Cells.Item(237, 11).Value = 500
$ws.Cells.Item(237, 12).Value = 500
$ws.Cells.Item(237, 13).Value = 500
$ws.Cells.Item(237, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(237, 15).Value = $origen
$ws.Cells.Item(237, 16).Value = 500
$ws.Cells.Item(237, 17).Value = 1
$ws.Cells.Item(237, 18).Value = $clasificacion
